$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: "Box" -> "Test"
$ws.Range("A2").Value = "Test"
$ws.Range("B2").Value = "Test"

# Remove row 3 entirely (A3/B3 "Piece")
$ws.Rows("3").Delete()

# Update the active selection to B3 (now an empty cell below the data)
$ws.Range("B3").Select()
